# fix missing edition in Mireya
# The "Mireya" sheet is missing the "No.9/Jun/1944" edition row between
# "No.7 y 8/Abr y May/1944" (row 8) and "No.10 y 11/Jul y Ago/1944" (row 9).
# Insert a new row at position 8, shifting the rest down, then fill it in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mireya")

# Insert a new blank row above current row 8 (pushes old row 8.. down by one)
$ws.Rows.Item(8).Insert()

# Column D holds the raw "No.X/Mon/Year" string; A/B/C are shared formulas
# parsing volume/month/year out of it. Replicate the same formulas used by
# the surrounding rows for the newly inserted row 8.
$ws.Range("A8").Formula = '=SUBSTITUTE(LEFT(D8,FIND("/",D8) - 1), "No.", "",1)'
$ws.Range("B8").Formula = '=MID(D8, FIND("/", D8) + 1, FIND( "!", SUBSTITUTE(D8, "/", "!", 2) ) - 2 - FIND("/", D8) + 1 )'
$ws.Range("C8").Formula = '=RIGHT(D8, 4)'
$ws.Range("D8").Value = "No.9/Jun/1944"

# Make Mireya the active sheet, matching the author's focus when fixing
# this row, and select the newly edited cell.
$ws.Activate()
$ws.Range("D9").Select()
